# Trey power network workbook update
# "Calcul power for each cabinet":
#  - bus 0 (STMT003438) becomes the 0.4 kV / PQ secondary bus
#  - a new bus 12 (STMT003438HV) is added as the 18 kV / Slack primary bus
#  - trafo 0 now links hv_bus 12 -> lv_bus 0 with its electrical characteristics
#  - ext_grid 0 is now attached to the new HV bus (12) instead of bus 0
#  - load table is filled in with one load per LV cabinet (name, bus, p_mw)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# bus sheet
# ---------------------------------------------------------------------------
$busWs = $wb.Worksheets.Item("bus")

# bus 0 (row 2) moves to the LV / PQ side
$busWs.Range("C2").Value = 0.4
$busWs.Range("D2").Value = "PQ"

# new bus 12 (row 14): the former HV / Slack side of the original bus
$busWs.Range("B14").Value = "STMT003438HV"
$busWs.Range("C14").Value = 18
$busWs.Range("D14").Value = "Slack"
$busWs.Range("E14").Value = "Trafo"

$busWs.Columns.Item(2).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# trafo sheet
# ---------------------------------------------------------------------------
$trafoWs = $wb.Worksheets.Item("trafo")

$trafoWs.Range("B2").Value = "STMT003438"
$trafoWs.Range("C2").Value = 12
$trafoWs.Range("D2").Value = 0
$trafoWs.Range("E2").Value = [double]"0.63"
$trafoWs.Range("F2").Value = [double]"18.3"
$trafoWs.Range("G2").Value = [double]"0.42"
$trafoWs.Range("H2").Value = 4
$trafoWs.Range("I2").Value = [double]"0.42"
$trafoWs.Range("J2").Value = [double]"0.65"
$trafoWs.Range("K2").Value = [double]"1.8"

# ---------------------------------------------------------------------------
# ext_grid sheet
# ---------------------------------------------------------------------------
$extGridWs = $wb.Worksheets.Item("ext_grid")
$extGridWs.Range("C2").Value = 12

# ---------------------------------------------------------------------------
# load sheet: one row per LV cabinet (name / bus idx / p_mw)
# ---------------------------------------------------------------------------
$loadWs = $wb.Worksheets.Item("load")

$loadWs.Range("B2").Value = "STMT003438"
$loadWs.Range("C2").Value = 0
$loadWs.Range("D2").Value = [double]"2.2329000000000002E-2"

$loadWs.Range("B3").Value = "CDBT004764"
$loadWs.Range("C3").Value = 1
$loadWs.Range("D3").Value = [double]"4.1394E-2"

$loadWs.Range("B4").Value = "CDBT003746"
$loadWs.Range("C4").Value = 2
$loadWs.Range("D4").Value = [double]"9.188E-3"

$loadWs.Range("B5").Value = "CDBT004760"
$loadWs.Range("C5").Value = 3
$loadWs.Range("D5").Value = [double]"2.2412999999999999E-2"

$loadWs.Range("B6").Value = "CDBT012139"
$loadWs.Range("C6").Value = 4
$loadWs.Range("D6").Value = [double]"6.8910000000000004E-3"

$loadWs.Range("B7").Value = "CDBT900784"
$loadWs.Range("C7").Value = 5
$loadWs.Range("D7").Value = [double]"6.352E-3"

$loadWs.Range("B8").Value = "CDBT901452"
$loadWs.Range("C8").Value = 6
$loadWs.Range("D8").Value = [double]"4.3998000000000002E-2"

$loadWs.Range("B9").Value = "CDBT004774"
$loadWs.Range("C9").Value = 7
$loadWs.Range("D9").Value = [double]"1.0565E-2"

$loadWs.Range("B10").Value = "CDBT901604"
$loadWs.Range("C10").Value = 8
$loadWs.Range("D10").Value = [double]"6.4539999999999997E-3"

$loadWs.Range("B11").Value = "CDBT016055"
$loadWs.Range("C11").Value = 9
$loadWs.Range("D11").Value = [double]"2.5128999999999999E-2"

$loadWs.Range("B12").Value = "N1"
$loadWs.Range("C12").Value = 10
$loadWs.Range("D12").Value = [double]"9.5980000000000006E-3"

# "60437" looks numeric, so it is entered with a leading apostrophe just like
# the matching cell on the bus sheet, keeping it as quote-prefixed text.
$loadWs.Range("B13").Value = "'60437"
$loadWs.Range("C13").Value = 11
$loadWs.Range("D13").Value = [double]"1.0946000000000001E-2"

# ---------------------------------------------------------------------------
# sheet selections / active sheet, matching the saved view state
# ---------------------------------------------------------------------------
$lineWs = $wb.Worksheets.Item("line")
$lineWs.Range("C2").Select() | Out-Null

$trafoWs.Range("I34").Select() | Out-Null

$extGridWs.Range("E36").Select() | Out-Null

$busWs.Range("A2:A14").Select() | Out-Null

$loadWs.Activate() | Out-Null
$loadWs.Range("G8").Select() | Out-Null
